$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 already carries the date-row formatting used by the log (fill, border,
# date number format) - copy it down onto B10 before filling in the value.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122) # xlPasteFormats

# Fill in the data for row 10 (new day of thesis-writing log)
$ws.Range("B10").Value = 45142
$ws.Range("C10").Formula = "=D10-D9"
$ws.Range("D10").Value = 1026

# Match the author's last cell selection after entering the data
$ws.Range("E12").Select()
